$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number format to text ("@") on column D for all rows that need it,
# so numeric-looking strings (e.g. "1.007", "30.106.04") are stored as text,
# matching the original inline-string cell type. Style is reset to Normal
# afterwards so no stray formatting is left on the cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '30.106.04'
$ws.Range('E2').Value = '  +0.10%  '

$ws.Range('D3').Value = '2.117.75'
$ws.Range('E3').Value = '  +0.57%  '

$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').Value = '346.51'
$ws.Range('E5').Value = '  +0.40%  '

$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  -0.07%  '

$ws.Range('D7').Value = '0.5199'
$ws.Range('E7').Value = '  +0.35%  '

$ws.Range('D8').Value = '0.4471'
$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('D9').Value = '54.08'
$ws.Range('E9').Value = '  +3.57%  '

$ws.Range('D10').Value = '0.09369'
$ws.Range('E10').Value = '  -1.63%  '

$ws.Range('D11').Value = '1.186'
$ws.Range('E11').Value = '  +0.71%  '

$ws.Range('D12').Value = '25.33'
$ws.Range('E12').Value = '  +0.12%  '

$ws.Range('D13').Value = '8.664'
$ws.Range('E13').Value = '  +6.97%  '

$ws.Range('D14').Value = '6.971'
$ws.Range('E14').Value = '  +3.15%  '

$ws.Range('D15').Value = '2.087.35'
$ws.Range('E15').Value = '  -1.00%  '

$ws.Range('D16').Value = '102.58'
$ws.Range('E16').Value = '  +3.02%  '

$ws.Range('E17').Value = '  -0.23%  '

$ws.Range('E18').Value = '  -0.11%  '

$ws.Range('E19').Value = '  +4.42%  '

$ws.Range('D20').Value = '0.06705'
$ws.Range('E20').Value = '  +0.06%  '

$ws.Range('D21').Value = '6.312'
$ws.Range('E21').Value = '  +1.94%  '

$ws.Range('D22').Value = '1.007'
$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').Value = '30.143.04'
$ws.Range('E23').Value = '  -0.03%  '

$ws.Range('D24').Value = '12.76'
$ws.Range('E24').Value = '  +0.42%  '

$ws.Range('D25').Value = '2.337'
$ws.Range('E25').Value = '  +0.63%  '

$ws.Range('D26').Value = '2.366.97'
$ws.Range('E26').Value = '  +0.47%  '

$ws.Range('D27').Value = '22.17'
$ws.Range('E27').Value = '  +0.49%  '

$ws.Range('D28').Value = '2.551'
$ws.Range('E28').Value = '  +0.56%  '

$ws.Range('D29').Value = '163.31'
$ws.Range('E29').Value = '  -0.56%  '

$ws.Range('D30').Value = '134.12'
$ws.Range('E30').Value = '  +0.35%  '

$ws.Range('D31').Value = '1.156'
$ws.Range('E31').Value = '  -0.43%  '

$ws.Range('D32').Value = '1.789'
$ws.Range('E32').Value = '  +9.69%  '

$ws.Range('E33').Value = '  +0.28%  '

$ws.Range('D34').Value = '6.290'
$ws.Range('E34').Value = '  +0.48%  '

$ws.Range('D35').Value = '6.779'
$ws.Range('E35').Value = '  +9.66%  '

$ws.Range('D36').Value = '3.968'
$ws.Range('E36').Value = '  +0.62%  '

$ws.Range('D37').Value = '10.79'
$ws.Range('E37').Value = '  +6.03%  '

$ws.Range('D38').Value = '0.02646'
$ws.Range('E38').Value = '  +2.77%  '

$ws.Range('D39').Value = '0.06878'
$ws.Range('E39').Value = '  +1.40%  '

$ws.Range('D40').Value = '0.7140'
$ws.Range('E40').Value = '  +2.57%  '

$ws.Range('D41').Value = '12.73'
$ws.Range('E41').Value = '  +1.93%  '

$ws.Range('D42').Value = '0.2248'

$ws.Range('D43').Value = '1.332'
$ws.Range('E43').Value = '  +1.63%  '

$ws.Range('D44').Value = '0.6986'
$ws.Range('E44').Value = '  +4.12%  '

$ws.Range('D45').Value = '14.75'
$ws.Range('E45').Value = '  +3.16%  '

$ws.Range('D46').Value = '2.404'
$ws.Range('E46').Value = '  +5.16%  '

$ws.Range('D47').Value = '1.007'
$ws.Range('E47').Value = '  +0.07%  '

$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '1.263'
$ws.Range('E48').Value = '  +7.63%  '

$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D49').Value = '3.633'
$ws.Range('E49').Value = '  -0.22%  '

$ws.Range('D50').Value = '0.00000000347'
$ws.Range('E50').Value = '  -2.73%  '

$ws.Range('D51').Value = '1.228'
$ws.Range('E51').Value = '  +0.29%  '

# Reset column D style back to Normal (no explicit style) to avoid leaving
# formatting artifacts beyond the text content changes.
$dRange.Style = "Normal"